# Refresh the "想去人数" (want-to-go count) values in column F across all
# four sheets (展览, 演出, 本地生活, 全部类型) to match a newer scrape snapshot.
# 全部类型 aggregates rows from the other three sheets, so its F values are
# updated independently to stay in sync.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 45
$ws.Range("F4").Value = 8487
$ws.Range("F5").Value = 8487
$ws.Range("F6").Value = 557
$ws.Range("F7").Value = 7498
$ws.Range("F8").Value = 1156
$ws.Range("F9").Value = 622
$ws.Range("F10").Value = 520
$ws.Range("F13").Value = 237
$ws.Range("F15").Value = 177
$ws.Range("F16").Value = 12418
$ws.Range("F18").Value = 18
$ws.Range("F19").Value = 2608
$ws.Range("F20").Value = 3757
$ws.Range("F23").Value = 3019
$ws.Range("F24").Value = 10
$ws.Range("F25").Value = 128
$ws.Range("F27").Value = 18
$ws.Range("F29").Value = 3402
$ws.Range("F30").Value = 81
$ws.Range("F31").Value = 348
$ws.Range("F32").Value = 1759
$ws.Range("F34").Value = 148
$ws.Range("F35").Value = 6138
$ws.Range("F38").Value = 1875
$ws.Range("F39").Value = 1267
$ws.Range("F40").Value = 47
$ws.Range("F41").Value = 936
$ws.Range("F42").Value = 5
$ws.Range("F43").Value = 180
$ws.Range("F45").Value = 201
$ws.Range("F46").Value = 1128
$ws.Range("F47").Value = 1116
$ws.Range("F48").Value = 1622
$ws.Range("F49").Value = 31
$ws.Range("F50").Value = 123

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 5
$ws.Range("F8").Value = 261
$ws.Range("F10").Value = 55
$ws.Range("F11").Value = 215
$ws.Range("F22").Value = 79
$ws.Range("F28").Value = 7

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 353
$ws.Range("F3").Value = 504

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 45
$ws.Range("F4").Value = 5
$ws.Range("F6").Value = 353
$ws.Range("F7").Value = 504
$ws.Range("F9").Value = 8487
$ws.Range("F11").Value = 557
$ws.Range("F12").Value = 7498
$ws.Range("F13").Value = 7498
$ws.Range("F14").Value = 622
$ws.Range("F15").Value = 520
$ws.Range("F16").Value = 237
$ws.Range("F17").Value = 261
$ws.Range("F18").Value = 177
$ws.Range("F20").Value = 12418
$ws.Range("F22").Value = 18
$ws.Range("F23").Value = 2608
$ws.Range("F24").Value = 2608
$ws.Range("F25").Value = 3757
$ws.Range("F26").Value = 10
$ws.Range("F27").Value = 128
$ws.Range("F29").Value = 18
$ws.Range("F32").Value = 3402
$ws.Range("F33").Value = 348
$ws.Range("F34").Value = 1759
$ws.Range("F36").Value = 148
$ws.Range("F37").Value = 6138
$ws.Range("F38").Value = 79
$ws.Range("F41").Value = 1875
$ws.Range("F43").Value = 1267
$ws.Range("F44").Value = 47
$ws.Range("F45").Value = 936
$ws.Range("F46").Value = 180
$ws.Range("F47").Value = 201
$ws.Range("F48").Value = 1128
$ws.Range("F49").Value = 1116
$ws.Range("F50").Value = 1622
$ws.Range("F51").Value = 31
$ws.Range("F52").Value = 123
